# Updated: st 16. 12. 2020
# A new day's AgTests/AgPosit (columns H/I) values were inserted at the top
# of the existing block (row 220), pushing the previously-entered values
# down by one row. The value that used to sit in the last row (285) now
# lands in the newly-extended last row (286), and the first row of the
# block (220) is left without H/I values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 220
$lastRow  = 286

# Capture the "before" values of columns H (8) and I (9) for the block,
# from the last row down to the first row, before we start overwriting.
$hVals = @{}
$iVals = @{}
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $hVals[$r] = $ws.Cells.Item($r, 8).Value2
    $iVals[$r] = $ws.Cells.Item($r, 9).Value2
}

# Shift every row's H/I value down into the next row, working from the
# bottom up so we never overwrite a value before it has been read.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $srcRow = $r - 1
    $ws.Cells.Item($r, 8).Value2 = $hVals[$srcRow]
    $ws.Cells.Item($r, 9).Value2 = $iVals[$srcRow]
}

# The first row of the block no longer has AgTests/AgPosit data.
$ws.Range("H220:I220").ClearContents()
